$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the description text for the "send_command" row (C5):
# "send command to ec" -> "send command to board"
$ws.Range("C5").Value = "send command to board"

# Move the active selection from B6 to E5
$ws.Range("E5").Select()

$wb.Save()
